# Apply updates described by the commit diff.
# Each "find" string below is unique within the document, so a simple
# Find/Replace (wdReplaceAll = 2) for each pair, executed in the same
# top-to-bottom order as they appear in the document, reproduces the
# target content exactly (including cases where a later "replace"
# value re-creates text that existed earlier in the document, and
# where the same replacement text is used for two different cells).

$d = $word.ActiveDocument

$replacements = @(
    @("2024-12-16 Monday", "2024-12-17 Tuesday"),
    @("52÷4=13, 0", "88÷5=17, 3"),
    @("35÷2=17, 1", "37÷2=18, 1"),
    @("95÷5=19, 0", "88÷4=22, 0"),
    @("94÷5=18, 4", "56÷3=18, 2"),
    @("89÷8=11, 1", "56÷3=18, 2"),
    @("92÷4=23, 0", "70÷4=17, 2"),
    @("91÷6=15, 1", "48÷6=8, 0"),
    @("31÷2=15, 1", "93÷5=18, 3"),
    @("57÷6=9, 3", "62÷9=6, 8"),
    @("10÷9=1, 1", "14÷6=2, 2"),
    @("54÷5=10, 4", "79÷8=9, 7"),
    @("85÷7=12, 1", "27÷8=3, 3"),
    @("54÷7=7, 5", "26÷2=13, 0"),
    @("56÷9=6, 2", "36÷3=12, 0"),
    @("55÷2=27, 1", "47÷3=15, 2"),
    @("62÷7=8, 6", "94÷5=18, 4"),
    @("74÷8=9, 2", "20÷9=2, 2"),
    @("51÷3=17, 0", "91÷4=22, 3"),
    @("28÷2=14, 0", "40÷9=4, 4"),
    @("80÷2=40, 0", "96÷9=10, 6"),
    @("69÷8=8, 5", "25÷5=5, 0"),
    @("71÷8=8, 7", "35÷3=11, 2"),
    @("54÷9=6, 0", "61÷9=6, 7"),
    @("88÷8=11, 0", "48÷5=9, 3"),
    @("16÷7=2, 2", "88÷3=29, 1")
)

foreach ($pair in $replacements) {
    $findText = $pair[0]
    $replaceText = $pair[1]

    $range = $d.Content
    $range.Find.Execute($findText, $true, $true, $false, $false, $false, `
                         $true, 1, $false, $replaceText, 2)
}
